$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (shifts existing D:K to F:M)
$ws.Range("D1:E1").EntireColumn.Insert()

# Copy cell formatting from column F into new columns D and E for all data rows
$ws.Range("F5:F102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$ws.Range("F5:F102").Copy()
$ws.Range("E5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new D and E columns with the latest two quarters of data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 45400
$ws.Range("E8").Value = 37400
$ws.Range("D9").Value = 21100
$ws.Range("E9").Value = 20400
$ws.Range("D10").Value = 24300
$ws.Range("E10").Value = 17000
$ws.Range("D12").Value = 16600
$ws.Range("E12").Value = 15800
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 71300
$ws.Range("E17").Value = 67500
$ws.Range("D18").Value = -25900
$ws.Range("E18").Value = -30100
$ws.Range("D20").Value = -4700
$ws.Range("E20").Value = 200
$ws.Range("D21").Value = -27300
$ws.Range("E21").Value = -26500
$ws.Range("D22").Value = 2100
$ws.Range("E22").Value = 1800
$ws.Range("D23").Value = -32600
$ws.Range("E23").Value = -31700
$ws.Range("D24").Value = -2800
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -29800
$ws.Range("E26").Value = -31700
$ws.Range("D27").Value = -29800
$ws.Range("E27").Value = -31700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 4700
$ws.Range("E32").Value = -200
$ws.Range("D33").Value = -29800
$ws.Range("E33").Value = -31700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -29800
$ws.Range("E35").Value = -31700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 112200
$ws.Range("E41").Value = 101400
$ws.Range("D42").Value = 13700
$ws.Range("E42").Value = 27800
$ws.Range("D43").Value = 26300
$ws.Range("E43").Value = 25500
$ws.Range("D44").Value = 0
$ws.Range("E44").Value = 0
$ws.Range("D45").Value = 13300
$ws.Range("E45").Value = 12700
$ws.Range("D46").Value = 165400
$ws.Range("E46").Value = 167300
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 400
$ws.Range("D48").Value = 27900
$ws.Range("E48").Value = 29300
$ws.Range("D49").Value = 80600
$ws.Range("E49").Value = 79000
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 9100
$ws.Range("E52").Value = 8500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 283000
$ws.Range("E54").Value = 284400
$ws.Range("D57").Value = 7800
$ws.Range("E57").Value = 8000
$ws.Range("D58").Value = 1900
$ws.Range("E58").Value = 10000
$ws.Range("D59").Value = 26600
$ws.Range("E59").Value = 24800
$ws.Range("D60").Value = 36300
$ws.Range("E60").Value = 42900
$ws.Range("D61").Value = 75900
$ws.Range("E61").Value = 52300
$ws.Range("D62").Value = 9000
$ws.Range("E62").Value = 9900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 121100
$ws.Range("E66").Value = 105000
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -516700
$ws.Range("E72").Value = -486900
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 161800
$ws.Range("E76").Value = 179400
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = -29800
$ws.Range("E81").Value = -31700
$ws.Range("D83").Value = 3300
$ws.Range("E83").Value = 3300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -15500
$ws.Range("E89").Value = -18100
$ws.Range("D91").Value = -1700
$ws.Range("E91").Value = -1200
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 11200
$ws.Range("E94").Value = 6500
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 16000
$ws.Range("E100").Value = 64100
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 11700
$ws.Range("E102").Value = 52600
